$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.699.97'
$ws.Range('E2').Value = '  -0.28%  '
$ws.Range('D3').Value = '1.847.64'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.014'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -2.48%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '319.25'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.70%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.010'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.75%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4311'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -2.61%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3742'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.61%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07344'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.69%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.8795'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.68%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '21.56'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.75%  '
$ws.Range('D12').Value = '1.837.97'
$ws.Range('E12').Value = '  -1.80%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.722'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.68%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.450'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.08%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.07135'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.46%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '87.89'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +4.86%  '
$ws.Range('E17').Value = '  -2.57%  '
$ws.Range('E18').Value = '  -1.91%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.010'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.75%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '15.46'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.73%  '
$ws.Range('D21').Value = '27.707.84'
$ws.Range('E21').Value = '  -0.28%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.243'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.47%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '11.14'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.93%  '
$ws.Range('D24').Value = '2.077.15'
$ws.Range('E24').Value = '  -1.14%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.009'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.52%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '155.55'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.26%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.60'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.54%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.131'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +7.49%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.377'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.81%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '120.23'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.99%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08929'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.57%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.226'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.99%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.7787'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.560'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.37%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.914'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -6.11%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.011'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.67%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.138'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.52%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.05328'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.46%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01971'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.18%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '7.228'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +4.69%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.878'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.46%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.5152'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.06%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1678'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.04%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.888'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.35%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '10.69'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.38%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '109.00'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.90%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.06515'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.60%  '
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.4728'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.31%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.697'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.63%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.010'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.72%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.876'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.32%  '
